# pptx: Fix list level numbering
#
# In PowerPoint, the content of a top-level list is at the same level as
# the content of a top-level paragraph. This presentation had every list
# paragraph pushed one level too deep (IndentLevel 2, i.e. OOXML lvl="1").
# Bring every such paragraph back up to the top level (IndentLevel 1,
# i.e. OOXML lvl="0") without touching paragraphs that are already at
# the top level (titles, intro/lead-in lines, etc).

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                $paras = $tr.Paragraphs()
                for ($pi = 1; $pi -le $paras.Count; $pi++) {
                    $para = $paras.Paragraphs($pi, 1)
                    if ($para.IndentLevel -eq 2) {
                        $para.IndentLevel = 1
                    }
                }
            }
        }
    }
}
